$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

$ws.Range("E18").Value = "ICD30449358"
$ws.Range("E19").Value = "ICD30449367"
$ws.Range("E20").Value = "ICD30449383"
$ws.Range("E21").Value = "ICD30449417"
$ws.Range("E22").Value = "ICD30449427"
$ws.Range("E23").Value = "ICD30452783"
$ws.Range("E24").Value = "ICD30453627"
$ws.Range("E25").Value = "ICD30453761"
$ws.Range("E26").Value = "ICD30453807"
